$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new date column before column DT ---
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Insert a new column at DT, shifting DT:EX (and their data) right to DU:EY
$wsPrix.Range("DT:DT").Insert()

# New header cell (DT1) continues the "13-nov"/"14-nov" sequence with "15-nov"
$wsPrix.Range("DT1").Value = "15-nov"
$wsPrix.Range("DS1").Copy()
$wsPrix.Range("DT1").PasteSpecial(-4122)

# New data column (DT2:DT25) has no values yet -> placeholder "-"
$wsPrix.Range("DT2:DT25").Value = "-"

# --- Sheet "CO2": append a new row of data ---
$wsCO2 = $wb.Worksheets.Item("CO2")
# Force the date to be stored as plain text (matching the existing rows),
# not auto-converted to a date serial number.
$wsCO2.Range("A152").NumberFormat = "@"
$wsCO2.Range("A152").Value = "2025-11-13"
$wsCO2.Range("A152").ClearFormats()
$wsCO2.Range("B152").Value = 81.02
